$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Stash pristine formatting samples into far-away scratch cells so we can
# restore the original style after writing date-like strings further down
# (see below: assigning a date-shaped literal like "09.01.2025" auto-detects
# as a real date unless the cell is pre-formatted as text, which otherwise
# leaves a stray NumberFormat behind).
$ws.Range("C1").Copy()
$ws.Range("Z1").PasteSpecial(-4122)   # xlPasteFormats: header-row style (s=3)

# --- Remove the dropped attendance row (old "Шекшукв Филипп"). Excel shifts
# everything below up automatically; there is nothing below row 4 here.
$ws.Rows.Item(4).Delete()

# --- Widen column A for the longer name.
$ws.Columns.Item(1).ColumnWidth = 23.17

# --- New lesson dates across the header row. Force text (not an auto-parsed
# date) by pre-formatting as Text, then repaint the original header style.
$headerCells = "B1", "C1", "D1", "E1"
$headerDates = "09.01.2025", "10.01.2025", "13.01.2025", "12.01.2025"
for ($i = 0; $i -lt $headerCells.Length; $i++) {
    $cell = $ws.Range($headerCells[$i])
    $cell.NumberFormat = "@"
    $cell.Value = $headerDates[$i]
    $ws.Range("Z1").Copy()
    $cell.PasteSpecial(-4122)
}

# --- Row 2: renamed student, attendance reset for the new date columns.
$ws.Range("A2").Value = "Точкееееееее Точкеее"
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = "З"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""

# --- Row 3: corrected surname spelling, updated attendance marks.
$ws.Range("A3").Value = "Шекшуев Филипп"
$ws.Range("B3").Value = "+"
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = "+"
$ws.Range("E3").Value = ""

# --- Drop the scratch cell used for format-painting.
$ws.Range("Z1").Clear()
